# Update Excel file via API
# Adds the latest (2019/2020) reporting-year rows across the "By the Numbers"
# dashboard sheets, drops the retired "Items delivered to faculty offices"
# metric column on CollectionUseDelivery, and leaves the selection on the
# SocialMedia sheet (mirrors the author's final click position).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Rankings: one more data point appended at the bottom (already sorted
# ascending by year, so this is a new last row, not an inserted one).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Rankings")
$ws.Cells.Item(18, 1).Value = 2019
$ws.Cells.Item(18, 2).Value = 35

# ---------------------------------------------------------------------
# EngageLearn: new 2020 row inserted above the existing data.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("EngageLearn")
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).Value = 2020
$ws.Cells.Item(2, 2).Value = 1272
$ws.Cells.Item(2, 3).Value = 17510
$ws.Cells.Item(2, 4).Value = 13872

# ---------------------------------------------------------------------
# Collections: new 2020 row inserted above the existing data.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Collections")
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).Value = 2020
$ws.Cells.Item(2, 2).Value = 2868659
$ws.Cells.Item(2, 3).Value = 2850467
$ws.Cells.Item(2, 4).Value = 39627

# ---------------------------------------------------------------------
# CollectionUseDelivery: drop the retired "Items delivered to faculty
# offices" column (I), then insert the new 2020 row above the rest.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CollectionUseDelivery")
$ws.Columns.Item(9).Delete()
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).Value = 2020
$ws.Cells.Item(2, 2).Value = 60268
$ws.Cells.Item(2, 3).Value = 2971571
$ws.Cells.Item(2, 4).Value = 537016
$ws.Cells.Item(2, 5).Value = 6412039
$ws.Cells.Item(2, 6).Value = 1298547
$ws.Cells.Item(2, 7).Value = 15106
$ws.Cells.Item(2, 8).Value = 15421

# ---------------------------------------------------------------------
# Expenditures: new 2020 row inserted above the existing data.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Expenditures")
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).Value = 2020
$ws.Cells.Item(2, 2).Value = 12.66
$ws.Cells.Item(2, 3).Value = 83.24
$ws.Cells.Item(2, 4).Value = 4.0999999999999996

# ---------------------------------------------------------------------
# SpacesStaff: new 2020 row inserted above the existing data.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SpacesStaff")
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).Value = 2020
$ws.Cells.Item(2, 2).Value = 10
$ws.Cells.Item(2, 3).Value = 311555
$ws.Cells.Item(2, 4).Value = 162

# ---------------------------------------------------------------------
# Visitors: new 2020 row inserted above the existing data.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Visitors")
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).Value = 2020
$ws.Cells.Item(2, 2).Value = 759258
$ws.Cells.Item(2, 3).Value = 10894129

# ---------------------------------------------------------------------
# SocialMedia: new 2020 row inserted above the existing data.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SocialMedia")
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).Value = 2020
$ws.Cells.Item(2, 2).Value = 2538
$ws.Cells.Item(2, 3).Value = 585
$ws.Cells.Item(2, 4).Value = 1009

# ---------------------------------------------------------------------
# Restore each sheet's cursor position; the last Select() made here wins
# as the workbook's active sheet/tab, matching SocialMedia being the
# front-most tab in the saved file.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Rankings").Range("A19").Select()
$wb.Worksheets.Item("EngageLearn").Range("E2").Select()
$wb.Worksheets.Item("Collections").Range("C2").Select()
$wb.Worksheets.Item("CollectionUseDelivery").Range("H18").Select()
$wb.Worksheets.Item("Expenditures").Range("E2").Select()
$wb.Worksheets.Item("SpacesStaff").Range("B5").Select()
$wb.Worksheets.Item("Visitors").Range("C3").Select()
$wb.Worksheets.Item("SocialMedia").Range("C2").Select()
